$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 239. This shifts the existing rows
# 239-254 down to 240-255, preserving their data/formatting unchanged.
$ws.Rows(239).Insert()

# Populate the newly inserted row 239 with the new weekly price entry.
$ws.Range("A239").Value = 5
$ws.Range("B239").Value = "Macroferia Regional de Talca"
$ws.Range("C239").Value = "Maule"
$ws.Range("D239").Value = 44714
$ws.Range("E239").Value = 7
$ws.Range("F239").Value = 100112009
$ws.Range("G239").Value = "Acelga"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 500
$ws.Range("K239").Value = 3000
$ws.Range("L239").Value = 3000
$ws.Range("M239").Value = 3000
$ws.Range("N239").Value = "$/docena de atados (4 kilos)"
$ws.Range("O239").Value = "Región del Maule"
$ws.Range("P239").Value = 750
$ws.Range("Q239").Value = 4
$ws.Range("R239").Value = "Hortaliza"
